$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.173.69"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.904.72"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  +0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "306.22"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5239"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.95%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3768"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.74%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07251"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.54%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "21.17"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.17%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.9031"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.23%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08518"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +11.42%  "
$ws.Range("D13").Value = "1.911.67"
$ws.Range("E13").Value = "  +0.99%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "95.92"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.03%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.294"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.41%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.07%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000008638"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "27.209.38"
$ws.Range("E20").Value = "  +0.33%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.069"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "2.171.90"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("E23").Value = "  +0.56%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.436"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "147.50"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.14%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.295"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.751"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "18.26"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.85%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "114.99"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.821"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.56%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.917"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.03%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.09282"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.02%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.8059"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.99%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.05052"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.80%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.240"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.41%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.452"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +4.89%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.959"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.10%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.623"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.32%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.5707"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.12%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.01999"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.08%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.073"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.23%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "9.174"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.81%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.644"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.17%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "116.24"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.30%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.1519"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("E46").Value = "  +1.46%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("E50").Value = "  -0.10%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "64.25"
$c.Style = "Normal"
